$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96
$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "Balanço Geral"
$ws.Cells.Item($row, 3).Value = "Economia"
$ws.Cells.Item($row, 4).Value = "2025-04-09T12:57"
$ws.Cells.Item($row, 5).Value = "Positivo"
$ws.Cells.Item($row, 6).Value = "Mutirão Facilita Imposto de Renda tem consultoria gratuita no Centro da cidade. Foi hoje no Centro da cidade. Repórter *ao vivo*. Muita gente aproveitou a oportunidade. Evento é parceria entre o Conselho Regional de Contabilidade do RJ e a Prefeitura de Campos. Equipe de contadores atendendo. Entrevista com delegada do CRC, Fabiana Viana."
